$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4812.727
$ws.Range("J2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("N2").Value = -1226

# Hunk 1: ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 657.21875
$ws.Range("I33").Value = 128.84
$ws.Range("K33").Value = 128.84
$ws.Range("M33").Value = 100.16

# Hunk 2: ALC!row40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 1425
$ws.Range("J40").Value = 1640.4
$ws.Range("K40").Value = 1425
$ws.Range("L40").Value = 1640.4
$ws.Range("M40").Value = -1250
$ws.Range("N40").Value = -1990.4

# Hunk 3: ALC!row51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8061.0557
$ws.Range("J51").Value = 2654.4546
$ws.Range("L51").Value = 2654.4546
$ws.Range("N51").Value = -3622.4546

# Hunk 4: ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4419.4
$ws.Range("I76").Value = 4419.4
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4419.4
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4104.4
$ws.Range("N76").ClearContents()

# Hunk 5: ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4419.4
$ws.Range("I79").Value = 4419.4
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4419.4
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3327.4
$ws.Range("N79").ClearContents()

# Hunk 6: ALC!row94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 3187.375
$ws.Range("I94").Value = 3214.2856
$ws.Range("J94").Value = 2999
$ws.Range("K94").Value = 3214.2856
$ws.Range("L94").Value = 2999
$ws.Range("M94").Value = -2763.2856
$ws.Range("N94").Value = -3901

# Hunk 7: ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1961.6111
$ws.Range("I138").Value = 1785.2667
$ws.Range("J138").Value = 2087.5715
$ws.Range("K138").Value = 5355.800099999999
$ws.Range("L138").Value = 6262.7145
$ws.Range("M138").Value = -215.8000999999995
$ws.Range("N138").Value = -16542.7145

# Hunk 8: ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 38092.63
$ws.Range("I2").Value = 1218.8889
$ws.Range("J2").Value = 111840.11
$ws.Range("K2").Value = 1218.8889
$ws.Range("L2").Value = 111840.11
$ws.Range("M2").Value = -1105.8889
$ws.Range("N2").Value = -112066.11

# Hunk 9: ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 112416.89
$ws.Range("I45").Value = 167718.5
$ws.Range("J45").Value = 1813.6666
$ws.Range("K45").Value = 167718.5
$ws.Range("L45").Value = 1813.6666
$ws.Range("M45").Value = -167341.5
$ws.Range("N45").Value = -2567.6666

# Hunk 10: ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 57401.277
$ws.Range("I102").Value = 168581.67
$ws.Range("K102").Value = 168581.67
$ws.Range("M102").Value = -166959.67

# Hunk 11: ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 77008660
$ws.Range("I110").Value = 111234024
$ws.Range("K110").Value = 111234024
$ws.Range("M110").Value = -111231979

# Hunk 12: ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 38092.63
$ws.Range("I116").Value = 1218.8889
$ws.Range("J116").Value = 111840.11
$ws.Range("K116").Value = 1218.8889
$ws.Range("L116").Value = 111840.11
$ws.Range("M116").Value = 1075.1111
$ws.Range("N116").Value = -116428.11

# Hunk 13: ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1826.3158
$ws.Range("I122").Value = 1638.75
$ws.Range("J122").Value = 2826.6667
$ws.Range("K122").Value = 4916.25
$ws.Range("L122").Value = 8480.000100000001
$ws.Range("M122").Value = -2466.25
$ws.Range("N122").Value = -13380.0001

# Hunk 14: ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2705.0588
$ws.Range("I132").Value = 2681.5898
$ws.Range("J132").Value = 2781.3333
$ws.Range("K132").Value = 8044.769400000001
$ws.Range("L132").Value = 8343.999899999999
$ws.Range("M132").Value = -5514.769400000001
$ws.Range("N132").Value = -13403.9999

# Hunk 15: BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 38092.63
$ws.Range("I3").Value = 1218.8889
$ws.Range("J3").Value = 111840.11
$ws.Range("K3").Value = 1218.8889
$ws.Range("L3").Value = 111840.11
$ws.Range("M3").Value = -1104.8889
$ws.Range("N3").Value = -112068.11

# Hunk 16: BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 42429.52
$ws.Range("I86").Value = 66382.06
$ws.Range("J86").Value = 1710.2
$ws.Range("K86").Value = 66382.06
$ws.Range("L86").Value = 1710.2
$ws.Range("M86").Value = -65259.06
$ws.Range("N86").Value = -3956.2

# Hunk 17: BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 42429.52
$ws.Range("I89").Value = 66382.06
$ws.Range("J89").Value = 1710.2
$ws.Range("K89").Value = 331910.3
$ws.Range("L89").Value = 8551
$ws.Range("M89").Value = -326294.3
$ws.Range("N89").Value = -19783

# Hunk 18: BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 76924024
$ws.Range("I107").Value = 200000420
$ws.Range("J107").Value = 1271.25
$ws.Range("K107").Value = 200000420
$ws.Range("L107").Value = 1271.25
$ws.Range("M107").Value = -199998500
$ws.Range("N107").Value = -5111.25

# Hunk 19: BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2375.9697
$ws.Range("I134").Value = 2155.72
$ws.Range("J134").Value = 3064.25
$ws.Range("K134").Value = 6467.16
$ws.Range("L134").Value = 9192.75
$ws.Range("M134").Value = -3932.16
$ws.Range("N134").Value = -14262.75

# Hunk 20: BSM!row140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 33945
$ws.Range("J140").Value = 33945
$ws.Range("L140").Value = 33945
$ws.Range("N140").Value = -44305

# Hunk 21: CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10568.154
$ws.Range("I99").Value = 3551.5
$ws.Range("J99").Value = 21794.8
$ws.Range("K99").Value = 3551.5
$ws.Range("L99").Value = 21794.8
$ws.Range("M99").Value = -2053.5
$ws.Range("N99").Value = -24790.8

# Hunk 22: CRP!row105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1094.5
$ws.Range("I105").Value = 1022.8571
$ws.Range("K105").Value = 1022.8571
$ws.Range("M105").Value = 724.1429000000001

# Hunk 23: CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10568.154
$ws.Range("I126").Value = 3551.5
$ws.Range("J126").Value = 21794.8
$ws.Range("K126").Value = 10654.5
$ws.Range("L126").Value = 65384.39999999999
$ws.Range("M126").Value = -8184.5
$ws.Range("N126").Value = -70324.39999999999

# Hunk 24: CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3834.1
$ws.Range("I132").Value = 3669.3635
$ws.Range("J132").Value = 4287.125
$ws.Range("K132").Value = 11008.0905
$ws.Range("L132").Value = 12861.375
$ws.Range("M132").Value = -8478.0905
$ws.Range("N132").Value = -17921.375

# Hunk 25: CUL!row34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 956.4286
$ws.Range("I34").Value = 250
$ws.Range("K34").Value = 750
$ws.Range("M34").Value = -666

# Hunk 26: CUL!row37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 2517976.5
$ws.Range("J37").Value = 2517976.5
$ws.Range("L37").Value = 7553929.5
$ws.Range("N37").Value = -7554153.5

# Hunk 27: CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 782.71
$ws.Range("J131").Value = 798.6701
$ws.Range("L131").Value = 2396.0103
$ws.Range("N131").Value = -12476.0103

# Hunk 28: GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1451.2222
$ws.Range("I113").Value = 1030
$ws.Range("J113").Value = 1613.2307
$ws.Range("K113").Value = 1030
$ws.Range("L113").Value = 1613.2307
$ws.Range("M113").Value = 1140
$ws.Range("N113").Value = -5953.2307

# Hunk 29: LTW!row61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1584.3529
$ws.Range("I61").Value = 1651.4
$ws.Range("J61").Value = 1488.5714
$ws.Range("K61").Value = 1651.4
$ws.Range("L61").Value = 1488.5714
$ws.Range("M61").Value = -1449.4
$ws.Range("N61").Value = -1892.5714

# Hunk 30: LTW!row108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Hunk 31: LTW!row113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1584.3529
$ws.Range("I113").Value = 1651.4
$ws.Range("J113").Value = 1488.5714
$ws.Range("K113").Value = 1651.4
$ws.Range("L113").Value = 1488.5714
$ws.Range("M113").Value = 518.5999999999999
$ws.Range("N113").Value = -5828.5714

# Hunk 32: WVR!row30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 24754.75
$ws.Range("I30").Value = 35004.5
$ws.Range("J30").Value = 14505
$ws.Range("K30").Value = 35004.5
$ws.Range("L30").Value = 14505
$ws.Range("M30").Value = -34897.5
$ws.Range("N30").Value = -14719

# Hunk 33: WVR!row46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 34895
$ws.Range("J46").Value = 34895
$ws.Range("L46").Value = 34895
$ws.Range("N46").Value = -35357

# Hunk 34: WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 400368
$ws.Range("I81").Value = 333647
$ws.Range("K81").Value = 667294
$ws.Range("M81").Value = -666233

# Hunk 35: WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 400368
$ws.Range("I84").Value = 333647
$ws.Range("K84").Value = 3336470
$ws.Range("M84").Value = -3331166

# Hunk 36: WVR!row113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 641.9032
$ws.Range("J113").Value = 951.0909
$ws.Range("L113").Value = 2853.2727
$ws.Range("N113").Value = -7193.2727

# Hunk 37: WVR!row134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 34895
$ws.Range("J134").Value = 34895
$ws.Range("L134").Value = 104685
$ws.Range("N134").Value = -109755
